$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top for the index/formula label, pushing the
# existing numeric data rows (and the shared-string table) down by one.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = " k*floor(x)+b+c1*sin(pi/7*(x-1))+c2*sin(2*pi/7*(x-1))"
$ws.Range("B1").Value = " "

# New unformatted row gets the sheet's auto-fit height for the default font
# instead of the explicit 15pt used by the original data rows.
$ws.Rows.Item(1).RowHeight = 12.8

$ws.Range("L12").Select()
